# Apply row-permutation updates to rows 7-21 of the Artfynd sheet.
# (species-occurrence records were re-ordered/re-matched to GPS points,
#  and the per-species "Taxonsorteringsordning" (col B) lookup values were refreshed.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("A7").Value = 111936876
$ws.Range("B7").Value = 89557
$ws.Range("Q7").Value = 449317
$ws.Range("R7").Value = 7087521

# Row 8
$ws.Range("A8").Value = 111936872
$ws.Range("B8").Value = 89557
$ws.Range("E8").Value = 5432
$ws.Range("F8").Value = "Granticka"
$ws.Range("G8").Value = "Porodaedalea chrysoloma"
$ws.Range("H8").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q8").Value = 449151
$ws.Range("R8").Value = 7087531
$ws.Range("AC8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()

# Row 9
$ws.Range("A9").Value = 111936800
$ws.Range("B9").Value = 56430
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 100109
$ws.Range("F9").Value = "Tretåig hackspett"
$ws.Range("G9").Value = "Picoides tridactylus"
$ws.Range("H9").Value = "(Linnaeus, 1758)"
$ws.Range("Q9").Value = 449178
$ws.Range("R9").Value = 7087530
$ws.Range("AC9").Value = "ringhack äldre"
$ws.Range("K9").Value = ""
$ws.Range("L9").Value = ""
$ws.Range("M9").Value = ""
$ws.Range("N9").Value = ""

# Row 10
$ws.Range("A10").Value = 111936855
$ws.Range("B10").Value = 78700
$ws.Range("E10").Value = 2081
$ws.Range("F10").Value = "Skrovellav"
$ws.Range("G10").Value = "Lobaria scrobiculata"
$ws.Range("H10").Value = "(Scop.) DC."
$ws.Range("Q10").Value = 449178
$ws.Range("R10").Value = 7087510

# Row 11
$ws.Range("A11").Value = 111936802
$ws.Range("B11").Value = 56430
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 100109
$ws.Range("F11").Value = "Tretåig hackspett"
$ws.Range("G11").Value = "Picoides tridactylus"
$ws.Range("H11").Value = "(Linnaeus, 1758)"
$ws.Range("Q11").Value = 448825
$ws.Range("R11").Value = 7087650
$ws.Range("AC11").Value = "ringhack äldre"
$ws.Range("K11").Value = ""
$ws.Range("L11").Value = ""
$ws.Range("M11").Value = ""
$ws.Range("N11").Value = ""

# Row 12
$ws.Range("A12").Value = 111936873
$ws.Range("B12").Value = 89557
$ws.Range("E12").Value = 5432
$ws.Range("F12").Value = "Granticka"
$ws.Range("G12").Value = "Porodaedalea chrysoloma"
$ws.Range("H12").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q12").Value = 449171
$ws.Range("R12").Value = 7087508
$ws.Range("AC12").ClearContents()
$ws.Range("K12").ClearContents()
$ws.Range("L12").ClearContents()
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()

# Row 13
$ws.Range("A13").Value = 111936877
$ws.Range("B13").Value = 89557
$ws.Range("Q13").Value = 449302
$ws.Range("R13").Value = 7087518

# Row 14
$ws.Range("A14").Value = 111936786
$ws.Range("B14").Value = 89539
$ws.Range("E14").Value = 1202
$ws.Range("F14").Value = "Ullticka"
$ws.Range("G14").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H14").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q14").Value = 449055
$ws.Range("R14").Value = 7087729

# Row 15
$ws.Range("B15").Value = 78700

# Row 16
$ws.Range("A16").Value = 111936895
$ws.Range("B16").Value = 85836
$ws.Range("E16").Value = 510
$ws.Range("F16").Value = "Doftskinn"
$ws.Range("G16").Value = "Cystostereum murrayi"
$ws.Range("H16").Value = "(Berk. & M.A. Curtis.) Pouzar"
$ws.Range("Q16").Value = 448925
$ws.Range("R16").Value = 7087774

# Row 17
$ws.Range("A17").Value = 111936791
$ws.Range("B17").Value = 90221
$ws.Range("D17").Value = "LC"
$ws.Range("E17").Value = 3298
$ws.Range("F17").Value = "Trådticka"
$ws.Range("G17").Value = "Climacocystis borealis"
$ws.Range("H17").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("Q17").Value = 449280
$ws.Range("R17").Value = 7087437

# Row 18
$ws.Range("A18").Value = 111936801
$ws.Range("B18").Value = 56430
$ws.Range("E18").Value = 100109
$ws.Range("F18").Value = "Tretåig hackspett"
$ws.Range("G18").Value = "Picoides tridactylus"
$ws.Range("H18").Value = "(Linnaeus, 1758)"
$ws.Range("Q18").Value = 448837
$ws.Range("R18").Value = 7087667
$ws.Range("AC18").Value = "ringhack äldre"
$ws.Range("K18").Value = ""
$ws.Range("L18").Value = ""
$ws.Range("M18").Value = ""
$ws.Range("N18").Value = ""

# Row 19
$ws.Range("A19").Value = 111936874
$ws.Range("B19").Value = 89557
$ws.Range("E19").Value = 5432
$ws.Range("F19").Value = "Granticka"
$ws.Range("G19").Value = "Porodaedalea chrysoloma"
$ws.Range("H19").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q19").Value = 449309
$ws.Range("R19").Value = 7087424

# Row 20
$ws.Range("A20").Value = 111936790
$ws.Range("B20").Value = 90221
$ws.Range("D20").Value = "LC"
$ws.Range("E20").Value = 3298
$ws.Range("F20").Value = "Trådticka"
$ws.Range("G20").Value = "Climacocystis borealis"
$ws.Range("H20").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("Q20").Value = 448944
$ws.Range("R20").Value = 7087698
$ws.Range("AC20").ClearContents()
$ws.Range("K20").ClearContents()
$ws.Range("L20").ClearContents()
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()

# Row 21
$ws.Range("A21").Value = 111936879
$ws.Range("B21").Value = 89557
$ws.Range("E21").Value = 5432
$ws.Range("F21").Value = "Granticka"
$ws.Range("G21").Value = "Porodaedalea chrysoloma"
$ws.Range("H21").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q21").Value = 449281
$ws.Range("R21").Value = 7087552

